$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Provas")

# New "Prova 02" column (C) scores, rows 7-26
$ws.Range("C7").Formula = "=20+15+8+15+25"
$ws.Range("C8").Value = 0
$ws.Range("C9").Formula = "=20+0+10+25+0"
$ws.Range("C10").Formula = "=20+0+8+25+25"
$ws.Range("C11").Formula = "=20+0+5+20+25"
$ws.Range("C12").Formula = "=20+0+13+22+25"
$ws.Range("C13").Value = 0
$ws.Range("C14").Formula = "=20+0+8+22+25"
$ws.Range("C15").Formula = "=20+15+13+10+25"
$ws.Range("C16").Formula = "=20+15+11+25+25"
$ws.Range("C17").Formula = "=20+7+8+25+25"
$ws.Range("C18").Formula = "=20+0+5+20+25"
$ws.Range("C19").Formula = "=20+0+11+12+15"
$ws.Range("C20").Formula = "=20+15+11+25+25"
$ws.Range("C21").Formula = "=20+7+8+25+25"
$ws.Range("C22").Formula = "=20+0+10+25+25"
$ws.Range("C23").Formula = "=20+15+13+22+25"
$ws.Range("C24").Formula = "=20+7+11+25+25"
$ws.Range("C25").Value = 0
$ws.Range("C26").Value = 100

# "Listas" keeps a remembered selection at A13 but is no longer the active tab
$wsListas = $wb.Worksheets.Item("Listas")
$wsListas.Activate()
$wsListas.Range("A13").Select()

# Active sheet moves to "Provas", selection at C27
$ws.Activate()
$ws.Range("C27").Select()
